$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two fixed bug entries (rows 2 and 4) - "Mario an nam lan 2" /
# "tu bien nho lai roi bien lon" and "Mario lon chet(do lot ho)" / "khong
# xet lai size..." - leaving the remaining rows where they are (gaps stay).
$ws.Range("A2:B2").ClearContents()
$ws.Range("A4:B4").ClearContents()

# Update the active selection to A7 as recorded after the edit.
$ws.Range("A7").Select()
